$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.64493191242218
$ws.Range("B1").Value = 2.415802717208862
$ws.Range("C1").Value = 2.790245771408081
$ws.Range("D1").Value = 3.419345617294312
$ws.Range("E1").Value = 1.207271933555603
